$wb = $excel.ActiveWorkbook
Write-Output $wb.GetType()
$members = $wb | Get-Member
Write-Output $members
